$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.602.77'
$ws.Range('E2').Value = '  +2.96%  '
$ws.Range('D3').Value = '3.147.77'
$ws.Range('E3').Value = '  +2.48%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '579.01'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = '180.08'
$ws.Range('E6').Value = '  +6.57%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').Value = '3.148.25'
$ws.Range('E8').Value = '  +2.64%  '
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('E10').Value = '  +1.19%  '
$ws.Range('E11').Value = '  +2.34%  '
$ws.Range('D12').Value = '0.470'
$ws.Range('E12').Value = '  +0.53%  '
$ws.Range('D13').Value = '0.0000244'
$ws.Range('E13').Value = '  +2.20%  '
$ws.Range('D14').Value = '37.17'
$ws.Range('E14').Value = '  +4.49%  '
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').Value = '68.529.05'
$ws.Range('E15').Value = '  +2.92%  '
$ws.Range('B16').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Value = '3.673.68'
$ws.Range('E16').Value = '  +2.55%  '
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').Value = '7.16'
$ws.Range('E18').Value = '  +3.09%  '
$ws.Range('D19').Value = '3.145.77'
$ws.Range('E19').Value = '  +2.40%  '
$ws.Range('D20').Value = '16.47'
$ws.Range('E20').Value = '  -2.59%  '
$ws.Range('D21').Value = '490.40'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '0.701'
$ws.Range('E22').Value = '  +2.30%  '
$ws.Range('E23').Value = '  +1.54%  '
$ws.Range('D24').Value = '84.01'
$ws.Range('E24').Value = '  +1.67%  '
$ws.Range('D25').Value = '2.36'
$ws.Range('E25').Value = '  +7.01%  '
$ws.Range('D26').Value = '13.06'
$ws.Range('E26').Value = '  +3.25%  '
$ws.Range('D27').Value = '10.61'
$ws.Range('E27').Value = '  +4.82%  '
$ws.Range('E28').Value = '  -0.03%  '
$ws.Range('D29').Value = '8.10'
$ws.Range('E29').Value = '  +4.16%  '
$ws.Range('D30').Value = '2.38'
$ws.Range('E30').Value = '  +5.44%  '
$ws.Range('E31').Value = '  +1.68%  '
$ws.Range('D32').Value = '28.28'
$ws.Range('E32').Value = '  +2.79%  '
$ws.Range('D33').Value = '0.113'
$ws.Range('E33').Value = '  +1.44%  '
$ws.Range('D34').Value = '0.0₃0956'
$ws.Range('E34').Value = '  +5.61%  '
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('D36').Value = '5.75'
$ws.Range('E36').Value = '  +3.14%  '
$ws.Range('D37').Value = '48.48'
$ws.Range('E37').Value = '  +3.21%  '
$ws.Range('D38').Value = '0.957'
$ws.Range('E38').Value = '  +1.08%  '
$ws.Range('E39').Value = '  +8.70%  '
$ws.Range('E40').Value = '  +5.02%  '
$ws.Range('E41').Value = '  +2.85%  '
$ws.Range('D42').Value = '49.21'
$ws.Range('E42').Value = '  +0.15%  '
$ws.Range('D43').Value = '8.43'
$ws.Range('E43').Value = '  +1.91%  '
$ws.Range('D44').Value = '2.76'
$ws.Range('E44').Value = '  +9.78%  '
$ws.Range('D45').Value = '400.94'
$ws.Range('E45').Value = '  +9.55%  '
$ws.Range('B46').Value = 'InjectiveProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D46').Value = '27.90'
$ws.Range('E46').Value = '  +13.43%  '
$ws.Range('B47').Value = 'Maker'
$ws.Range('C47').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D47').Value = '2.810.18'
$ws.Range('E47').Value = '  +1.76%  '
$ws.Range('D48').Value = '0.0350'
$ws.Range('E48').Value = '  +1.82%  '
$ws.Range('D49').Value = '134.39'
$ws.Range('E49').Value = '  -0.69%  '
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('D51').Value = '2.39'
$ws.Range('E51').Value = '  +11.12%  '
